$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 57968.938
$ws.Range("J17").Value = 57968.938
$ws.Range("L17").Value = 173906.814
$ws.Range("N17").Value = -174242.814
$ws.Range("H70").Value = 893.6111
$ws.Range("I70").Value = 749
$ws.Range("J70").Value = 1399.75
$ws.Range("K70").Value = 2247
$ws.Range("L70").Value = 4199.25
$ws.Range("M70").Value = -1977
$ws.Range("N70").Value = -4739.25
$ws.Range("H73").Value = 893.6111
$ws.Range("I73").Value = 749
$ws.Range("J73").Value = 1399.75
$ws.Range("K73").Value = 2247
$ws.Range("L73").Value = 4199.25
$ws.Range("M73").Value = -1311
$ws.Range("N73").Value = -6071.25
$ws.Range("H99").Value = 73344620
$ws.Range("I99").Value = 16692.6
$ws.Range("K99").Value = 50077.8
$ws.Range("M99").Value = -48579.8
$ws.Range("H131").Value = 1198
$ws.Range("I131").Value = 1198
$ws.Range("K131").Value = 3594
$ws.Range("M131").Value = 1446
$ws.Range("H135").Value = 694.4783
$ws.Range("I135").Value = 587.5238
$ws.Range("J135").Value = 1817.5
$ws.Range("K135").Value = 5287.7142
$ws.Range("L135").Value = 16357.5
$ws.Range("M135").Value = -2752.7142
$ws.Range("N135").Value = -21427.5
$ws.Range("H137").Value = 3564.6155
$ws.Range("I137").Value = 1546.341
$ws.Range("K137").Value = 4639.022999999999
$ws.Range("M137").Value = -2089.022999999999
$ws.Range("H138").Value = 2227.7954
$ws.Range("I138").Value = 1545.9131
$ws.Range("J138").Value = 2974.6191
$ws.Range("K138").Value = 4637.7393
$ws.Range("L138").Value = 8923.8573
$ws.Range("M138").Value = 502.2606999999998
$ws.Range("N138").Value = -19203.8573

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 998.9259
$ws.Range("I61").Value = 846.36
$ws.Range("K61").Value = 846.36
$ws.Range("M61").Value = -634.36
$ws.Range("H74").Value = 1633
$ws.Range("I74").Value = 1492.75
$ws.Range("J74").Value = 4999
$ws.Range("K74").Value = 1492.75
$ws.Range("L74").Value = 4999
$ws.Range("M74").Value = -618.75
$ws.Range("N74").Value = -6747
$ws.Range("H77").Value = 1633
$ws.Range("I77").Value = 1492.75
$ws.Range("J77").Value = 4999
$ws.Range("K77").Value = 7463.75
$ws.Range("L77").Value = 24995
$ws.Range("M77").Value = -3095.75
$ws.Range("N77").Value = -33731
$ws.Range("H122").Value = 1493.5454
$ws.Range("I122").Value = 1477.9
$ws.Range("K122").Value = 4433.700000000001
$ws.Range("M122").Value = -1983.700000000001
$ws.Range("H131").Value = 79999
$ws.Range("J131").Value = 79999
$ws.Range("L131").Value = 79999
$ws.Range("N131").Value = -90079
$ws.Range("H132").Value = 4549.5557
$ws.Range("I132").Value = 2995
$ws.Range("K132").Value = 8985
$ws.Range("M132").Value = -6455
$ws.Range("H136").Value = 998.9259
$ws.Range("I136").Value = 846.36
$ws.Range("K136").Value = 2539.08
$ws.Range("M136").Value = 10.92000000000007

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H129").Value = 66637
$ws.Range("J129").Value = 66637
$ws.Range("L129").Value = 66637
$ws.Range("N129").Value = -76637

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1664
$ws.Range("I31").Value = 1341.421
$ws.Range("K31").Value = 1341.421
$ws.Range("M31").Value = -1046.421
$ws.Range("H34").Value = 1664
$ws.Range("I34").Value = 1341.421
$ws.Range("K34").Value = 1341.421
$ws.Range("M34").Value = -1139.421
$ws.Range("H39").Value = 5560.8335
$ws.Range("I39").Value = 6125.5
$ws.Range("J39").Value = 5278.5
$ws.Range("K39").Value = 6125.5
$ws.Range("L39").Value = 5278.5
$ws.Range("M39").Value = -5734.5
$ws.Range("N39").Value = -6060.5
$ws.Range("H49").Value = 5560.8335
$ws.Range("I49").Value = 6125.5
$ws.Range("J49").Value = 5278.5
$ws.Range("K49").Value = 6125.5
$ws.Range("L49").Value = 5278.5
$ws.Range("M49").Value = -5943.5
$ws.Range("N49").Value = -5642.5
$ws.Range("H58").Value = 2028.4166
$ws.Range("I58").Value = 2072.3333
$ws.Range("K58").Value = 2072.3333
$ws.Range("M58").Value = -1869.3333
$ws.Range("H123").Value = 100924.336
$ws.Range("J123").Value = 100924.336
$ws.Range("L123").Value = 100924.336
$ws.Range("N123").Value = -110724.336
$ws.Range("H136").Value = 2028.4166
$ws.Range("I136").Value = 2072.3333
$ws.Range("K136").Value = 6216.999899999999
$ws.Range("M136").Value = -3666.999899999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12
$ws.Range("J2").Value = 13
$ws.Range("L2").Value = 78
$ws.Range("N2").Value = -304
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H98").Value = 1339
$ws.Range("I98").Value = 728.8333
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 2186.4999
$ws.Range("L98").Value = 15000
$ws.Range("M98").Value = -688.4998999999998
$ws.Range("N98").Value = -17996
$ws.Range("H119").Value = 999.5
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H128").Value = 442003.5
$ws.Range("I128").Value = 442003.5
$ws.Range("K128").Value = 1326010.5
$ws.Range("M128").Value = -1321030.5
$ws.Range("H129").Value = 127669.375
$ws.Range("J129").Value = 3374.0908
$ws.Range("L129").Value = 10122.2724
$ws.Range("N129").Value = -20122.2724

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3404.4736
$ws.Range("I126").Value = 3405.25
$ws.Range("K126").Value = 10215.75
$ws.Range("M126").Value = -7745.75
$ws.Range("H132").Value = 1931.5807
$ws.Range("I132").Value = 1827.6
$ws.Range("K132").Value = 5482.799999999999
$ws.Range("M132").Value = -2952.799999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 691.2941
$ws.Range("I55").Value = 543.6
$ws.Range("K55").Value = 543.6
$ws.Range("M55").Value = -370.6
$ws.Range("H82").Value = 1226.4546
$ws.Range("I82").Value = 1249.1
$ws.Range("K82").Value = 1249.1
$ws.Range("M82").Value = -888.0999999999999
$ws.Range("H85").Value = 1226.4546
$ws.Range("I85").Value = 1249.1
$ws.Range("K85").Value = 1249.1
$ws.Range("M85").Value = -1.099999999999909
$ws.Range("H136").Value = 4659.846
$ws.Range("I136").Value = 4062.0908
$ws.Range("K136").Value = 12186.2724
$ws.Range("M136").Value = -9636.2724

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3392.3572
$ws.Range("I81").Value = 3369.8518
$ws.Range("K81").Value = 6739.7036
$ws.Range("M81").Value = -5678.7036
$ws.Range("H84").Value = 3392.3572
$ws.Range("I84").Value = 3369.8518
$ws.Range("K84").Value = 33698.518
$ws.Range("M84").Value = -28394.518
$ws.Range("H132").Value = 10361.23
$ws.Range("I132").Value = 7760.4443
$ws.Range("J132").Value = 16213
$ws.Range("K132").Value = 23281.3329
$ws.Range("L132").Value = 48639
$ws.Range("M132").Value = -20751.3329
$ws.Range("N132").Value = -53699
